$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New query text for each tab's 'query' column (B) and the shared 'StatQuery' (C)
$participantsQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.reference_genome_assembly in ['GRCh38'] 
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id limit 100
'@

$samplesQuery = @'
MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE g.reference_genome_assembly in ['GRCh38']
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(samp.sample_tumor_status,'') as `Tumor`,
    coalesce(samp.sample_type,'') as `Analyte Type`
ORDER BY samp.sample_id limit 100
'@

$filesQuery = @'
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE g.reference_genome_assembly in ['GRCh38']
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name,'') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id, '') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER BY f.file_name limit 100
'@

$statQuery = @'
CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.reference_genome_assembly in ['GRCh38']
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.reference_genome_assembly in ['GRCh38']
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.reference_genome_assembly in ['GRCh38']
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS `Files`
'@

# Row 2 = ParticipantsTab, Row 3 = SamplesTab, Row 4 = FilesTab
$ws.Range("B2").Value = $participantsQuery
$ws.Range("C2").Value = $statQuery

$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $statQuery

$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQuery

# Re-assert the original explicit row heights: the new (longer) StatQuery
# text in column C would otherwise cause the wrap-text rows to auto-fit to
# a taller height.
$ws.Rows.Item(2).RowHeight = 242.25
$ws.Rows.Item(3).RowHeight = 260.25
$ws.Rows.Item(4).RowHeight = 279.75

# Update the view: clear the frozen/scrolled top-left cell back to A1 and
# change the selection from a merged-looking D4:E4 range down to just D4.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("D4").Select()
